# Updated symbol list on Sun Jan  8 19:59:27 UTC 2023 with GitHub Actions
#
# Refreshes the crypto price/volume snapshot in the sheet: updated Price
# (col D) and Volume(1h) (col E) figures for a batch of coins, plus a
# ranking swap between CoinbaseStockToken and BOLO (rows 47-48).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data stores prices/percentages as literal text (not numbers).
# Force the Text number format before assigning so these numeric-looking
# strings ("267.44", "2.18%", ...) stay text instead of being coerced to
# numeric/percentage values.
$numericCells = @("D2", "E2", "D3", "E3", "D4", "E4", "E5", "D6", "E6", "E7", "D8", "E8", "D9", "E9", "D10", "D11", "E11", "D12", "E12", "D13", "E13", "E14", "D15", "E15", "D16", "E16", "D18", "E18", "D19", "E19", "E21", "D22", "E22", "E23", "E24", "D25", "E25", "E27", "D40", "E40", "D41", "E41", "D42", "E42", "D44", "E44", "D45", "E45", "E46", "D47", "E47", "D48", "E48", "E49", "E50")
foreach ($addr in $numericCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "267.44"
$ws.Range("E2").Value = "2.18%"

$ws.Range("D3").Value = "26.69"
$ws.Range("E3").Value = "-1.77%"

$ws.Range("D4").Value = "4.688"
$ws.Range("E4").Value = "-0.41%"

$ws.Range("E5").Value = "-1.86%"

$ws.Range("D6").Value = "6.739"
$ws.Range("E6").Value = "0.13%"

$ws.Range("E7").Value = "0.13%"

$ws.Range("D8").Value = "0.9064"
$ws.Range("E8").Value = "-1.41%"

$ws.Range("D9").Value = "0.1414"
$ws.Range("E9").Value = "0.16%"

$ws.Range("D10").Value = "0.04901"

$ws.Range("D11").Value = "0.07108"
$ws.Range("E11").Value = "0.39%"

$ws.Range("D12").Value = "0.03207"
$ws.Range("E12").Value = "2.44%"

$ws.Range("D13").Value = "0.09022"
$ws.Range("E13").Value = "-0.38%"

$ws.Range("E14").Value = "-0.53%"

$ws.Range("D15").Value = "0.0006047"
$ws.Range("E15").Value = "-1.99%"

$ws.Range("D16").Value = "0.006020"
$ws.Range("E16").Value = "0.14%"

$ws.Range("D18").Value = "3.170"
$ws.Range("E18").Value = "0.20%"

$ws.Range("D19").Value = "2.241"
$ws.Range("E19").Value = "3.67%"

$ws.Range("E21").Value = "-0.83%"

$ws.Range("D22").Value = "4.065"
$ws.Range("E22").Value = "-0.56%"

$ws.Range("E23").Value = "0.02%"

$ws.Range("E24").Value = "-2.74%"

$ws.Range("D25").Value = "0.004134"
$ws.Range("E25").Value = "8.73%"

$ws.Range("E27").Value = "5.01%"

$ws.Range("D40").Value = "0.03913"
$ws.Range("E40").Value = "-0.34%"

$ws.Range("D41").Value = "0.1114"
$ws.Range("E41").Value = "0.06%"

$ws.Range("D42").Value = "0.004214"
$ws.Range("E42").Value = "1.90%"

$ws.Range("D44").Value = "0.01267"
$ws.Range("E44").Value = "-8.18%"

$ws.Range("D45").Value = "0.00005133"
$ws.Range("E45").Value = "-0.64%"

$ws.Range("E46").Value = "0.00%"

# Rows 47/48 swap rank: BOLO moves up to 47, CoinbaseStockToken drops to 48.
$ws.Range("B47").Value = "BOLO"
$ws.Range("C47").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D47").Value = "0.1247"
$ws.Range("E47").Value = "-25.64%"

$ws.Range("B48").Value = "CoinbaseStockToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D48").Value = "0.02448"
$ws.Range("E48").Value = "-31.80%"

$ws.Range("E49").Value = "0.00%"
$ws.Range("E50").Value = "0.00%"
